$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(1.417691046532468, 1.330987364808379, 1.277840640828146, 1.256206196152789, 1.252615241069549, 1.27754877591849, 1.387777398406342, 1.604628083676971, 1.764364128440263, 1.837122934134641, 1.864688041249053, 1.858750836489605, 1.839390476333108, 1.827533362919553, 1.759611010136439, 1.717966501878948, 1.694022570718971, 1.685917112632978, 1.722398715420297, 1.845076737196507, 1.925329009452639, 1.882490179271429, 1.720394917687258, 1.545890885390918)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}

$colD = @(0.01380876848636348, 0.01233916959299819, 0.01143558936547606, 0.01106709253718208, 0.01100588791616275, 0.0114306207823347, 0.01330232636778561, 0.0169617340021162, 0.01964236346696424, 0.02085990210918709, 0.02132065662397054, 0.02122143872608007, 0.0208978148535266, 0.02069954585874001, 0.01956275388067752, 0.0188648632006263, 0.01846327840129902, 0.01832727925868483, 0.01893917333703143, 0.02099287942120753, 0.02233332829308665, 0.02161807639537727, 0.01890557883736932, 0.01597307603675802)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colD[$i]
}

$colE = @(1.389111181297864, 1.279515203907891, 1.21243476335448, 1.185155240567099, 1.180628987074158, 1.212066630157409, 1.351280246533378, 1.625844798539674, 1.828394740440842, 1.920695385099521, 1.955667831947153, 1.948135034667416, 1.923572200121356, 1.908529304257939, 1.822365650899997, 1.769545940529468, 1.739180627951953, 1.72890214894116, 1.775167137278004, 1.93078637353284, 2.032609034598835, 1.978254596100982, 1.772625789539489, 1.551413601501253)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $colE[$i]
}

$colF = @(0.4492132562603928, 0.4458500536640457, 0.4443827027580838, 0.443934017012154, 0.4438684988853794, 0.4443760485825194, 0.4479289643344089, 0.4596852700100698, 0.4713083121899331, 0.4772585141204502, 0.4796081325998642, 0.479097795141513, 0.4774498802686935, 0.4764530721663789, 0.4709328836584064, 0.4677169697748695, 0.465929583855285, 0.4653350763402102, 0.4680528512067497, 0.4779312876741528, 0.4849498518904483, 0.4811520884874625, 0.4679008077671085, 0.4559853305767021)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $colF[$i]
}

$colG = @(0.002375239994871183, 0.002378748244697512, 0.002381014528822858, 0.002381966361171865, 0.00238212612425661, 0.002381027250682409, 0.002376426404435525, 0.002368290399221683, 0.002362847442738553, 0.002360486156297314, 0.002359608405811746, 0.002359796716190076, 0.002360413614642344, 0.002360793618446978, 0.002363004060260755, 0.002364389423947933, 0.002365197053584125, 0.002365472361208263, 0.002364240832073372, 0.002360231971683722, 0.002357707604953633, 0.002359046181705636, 0.002364307975650058, 0.002370397117822135)
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $colG[$i]
}

$colL = @(0.5052342562552781, 0.4459129144951248, 0.4093542267308692, 0.3944229347165447, 0.3919416118537526, 0.4091529921663835, 0.4848086715875013, 0.632076779787667, 0.7395942045027937, 0.7883571304505494, 0.8068008533324189, 0.8028296430134674, 0.7898749454001575, 0.7819369719644271, 0.7364044202760169, 0.7084335425127222, 0.6923316012133967, 0.6868774024024162, 0.7114125259364812, 0.7936806483996008, 0.8473206743509252, 0.8187037751766013, 0.7100657932583658, 0.5923554272453941)
for ($i = 0; $i -lt $colL.Length; $i++) {
    $ws.Cells.Item($i + 2, 12).Value = $colL[$i]
}

$colM = @(0.4051812384533235, 0.3694234055495826, 0.3474247125586203, 0.3384495142387038, 0.336958561838415, 0.3473037121855569, 0.3928610883439063, 0.4818469686495064, 0.5470059473961442, 0.576600734221401, 0.5878007019829283, 0.5853889011070237, 0.5775223027940086, 0.5727028732462855, 0.5450709021988303, 0.5281075122017427, 0.518346270386985, 0.5150405432449929, 0.5299137460897612, 0.5798331037029953, 0.6124177499826402, 0.5950305049898361, 0.5290971750251998, 0.4578121307324565)
for ($i = 0; $i -lt $colM.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $colM[$i]
}

$colO = @(1.411580729877301, 1.415623111248635, 1.420147486774567, 1.422501170409078, 1.422922692770129, 1.420177169337194, 1.412548490253982, 1.413957237635515, 1.425199111847576, 1.432578714682364, 1.435702889118829, 1.435015317858046, 1.432829112892506, 1.431533048417464, 1.42476275975244, 1.421192452429494, 1.419351924756143, 1.418765221407909, 1.421550442933352, 1.433462275943185, 1.443170769478428, 1.437811887643733, 1.42138793506092, 1.411800159689363)
for ($i = 0; $i -lt $colO.Length; $i++) {
    $ws.Cells.Item($i + 2, 15).Value = $colO[$i]
}
